$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Inflammatory-Mac -> FAPs (unchanged pair), values updated
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.7255343333333334
$ws.Range("H2").Value = 2.176603
$ws.Range("I2").Value = 0.6805708332369251
$ws.Range("J2").Value = 0.680570833236925
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.001424707585888889
$ws.Range("R2").Value = 0.012822368273
$ws.Range("S2").Value = 0.6805708332369251
$ws.Range("T2").Value = 0.680570833236925

# Row 3: now MuSCs -> FAPs (was Inflammatory-Mac -> MuSCs)
$ws.Range("A3").Value = "MuSCs"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05044199999999999
$ws.Range("H3").Value = 0.151326
$ws.Range("I3").Value = 0.04731596065539324
$ws.Range("J3").Value = 0.04731596065539324
$ws.Range("M3").Value = 0.001963666666666667
$ws.Range("N3").Value = 0.005891
$ws.Range("O3").Value = 1
$ws.Range("P3").Value = 1
$ws.Range("Q3").Value = 0.000099051274
$ws.Range("R3").Value = 0.0008914614659999999
$ws.Range("S3").Value = 0.04731596065539324
$ws.Range("T3").Value = 0.04731596065539324

# Row 4: now Resolving-Mac -> FAPs (was MuSCs -> FAPs)
$ws.Range("A4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.290091
$ws.Range("H4").Value = 0.870273
$ws.Range("I4").Value = 0.2721132061076817
$ws.Range("J4").Value = 0.2721132061076817
$ws.Range("O4").Value = 1
$ws.Range("P4").Value = 1
$ws.Range("Q4").Value = 0.0005696420270000001
$ws.Range("R4").Value = 0.005126778243000001
$ws.Range("S4").Value = 0.2721132061076817
$ws.Range("T4").Value = 0.2721132061076817

# Remove rows 5-7 (old extra rows no longer present)
$ws.Range("A5:T7").Delete()
